# Update each of the 6 sheets: fill in B102 (previously placeholder 0)
# with the real value, and append a new row 103 with the next trading
# date (serial 45966) in column A (formatted like the existing date
# column) and its value in column B.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Index = 1; B102 = 513700; B103 = 518529 },
    @{ Index = 2; B102 = 372983; B103 = 370671 },
    @{ Index = 3; B102 = 170442; B103 = 164224 },
    @{ Index = 4; B102 = 249792; B103 = 242799 },
    @{ Index = 5; B102 = 595164; B103 = 573768 },
    @{ Index = 6; B102 = 62696;  B103 = 60869 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Index)

    # Fill in the real value for the existing last row (102).
    $ws.Cells.Item(102, 2).Value = $u.B102

    # Add the new row (103) with date + value, matching the date
    # formatting used by the rest of column A.
    $ws.Cells.Item(103, 1).Value = 45966
    $ws.Cells.Item(103, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item(103, 2).Value = $u.B103
}
